# Aragon hospital COVID-19 occupancy tracker: append the daily hospital-level
# bed/ICU occupancy rows published for 2020-09-06 (date serial 44080, sourced
# from the 2020-09-05 report) and 2020-09-07 (date serial 44081, sourced from
# the 2020-09-06 report) - 22 hospitals per day, 44 rows total.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (2922) into the new rows
$styleSrc = $ws.Range("A2922:H2922")
for ($r = 2923; $r -le 2966; $r++) {
    $styleSrc.Copy()
    $ws.Range("A" + $r + ":H" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Populate the new rows (2923:2966) with the updated hospital data for 2020-09-06 and 2020-09-07
$ws.Cells.Item(2923, 1).Value = 44080
$ws.Cells.Item(2923, 2).Value = "Hospital Universitario Miguel Servet"
$ws.Cells.Item(2923, 3).Value = 99
$ws.Cells.Item(2923, 4).Value = 23
$ws.Cells.Item(2923, 5).Value = "Zaragoza"
$ws.Cells.Item(2923, 6).Value = "Zaragoza"
$ws.Cells.Item(2923, 7).Value = 50297
$ws.Cells.Item(2923, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2924, 1).Value = 44080
$ws.Cells.Item(2924, 2).Value = "Hospital Clínico Universitario"
$ws.Cells.Item(2924, 3).Value = 144
$ws.Cells.Item(2924, 4).Value = 22
$ws.Cells.Item(2924, 5).Value = "Zaragoza"
$ws.Cells.Item(2924, 6).Value = "Zaragoza"
$ws.Cells.Item(2924, 7).Value = 50297
$ws.Cells.Item(2924, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2925, 1).Value = 44080
$ws.Cells.Item(2925, 2).Value = "Hospital Royo Villanova"
$ws.Cells.Item(2925, 3).Value = 42
$ws.Cells.Item(2925, 4).Value = 4
$ws.Cells.Item(2925, 5).Value = "Zaragoza"
$ws.Cells.Item(2925, 6).Value = "Zaragoza"
$ws.Cells.Item(2925, 7).Value = 50297
$ws.Cells.Item(2925, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2926, 1).Value = 44080
$ws.Cells.Item(2926, 2).Value = "Hospital Nuestra Señora de Gracia"
$ws.Cells.Item(2926, 3).Value = 11
$ws.Cells.Item(2926, 4).Value = 0
$ws.Cells.Item(2926, 5).Value = "Zaragoza"
$ws.Cells.Item(2926, 6).Value = "Zaragoza"
$ws.Cells.Item(2926, 7).Value = 50297
$ws.Cells.Item(2926, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2927, 1).Value = 44080
$ws.Cells.Item(2927, 2).Value = "Hospital General de la Defensa"
$ws.Cells.Item(2927, 3).Value = 18
$ws.Cells.Item(2927, 4).Value = 2
$ws.Cells.Item(2927, 5).Value = "Zaragoza"
$ws.Cells.Item(2927, 6).Value = "Zaragoza"
$ws.Cells.Item(2927, 7).Value = 50297
$ws.Cells.Item(2927, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2928, 1).Value = 44080
$ws.Cells.Item(2928, 2).Value = "Hospital Obispo Polanco"
$ws.Cells.Item(2928, 3).Value = 21
$ws.Cells.Item(2928, 4).Value = 2
$ws.Cells.Item(2928, 5).Value = "Teruel"
$ws.Cells.Item(2928, 6).Value = "Teruel"
$ws.Cells.Item(2928, 7).Value = 44216
$ws.Cells.Item(2928, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2929, 1).Value = 44080
$ws.Cells.Item(2929, 2).Value = "Hospital de Alcañiz"
$ws.Cells.Item(2929, 3).Value = 19
$ws.Cells.Item(2929, 4).Value = 0
$ws.Cells.Item(2929, 5).Value = "Alcañiz"
$ws.Cells.Item(2929, 6).Value = "Teruel"
$ws.Cells.Item(2929, 7).Value = 44216
$ws.Cells.Item(2929, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2930, 1).Value = 44080
$ws.Cells.Item(2930, 2).Value = "Hospital de Barbastro"
$ws.Cells.Item(2930, 3).Value = 13
$ws.Cells.Item(2930, 4).Value = 3
$ws.Cells.Item(2930, 5).Value = "Barbastro"
$ws.Cells.Item(2930, 6).Value = "Huesca"
$ws.Cells.Item(2930, 7).Value = 22125
$ws.Cells.Item(2930, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2931, 1).Value = 44080
$ws.Cells.Item(2931, 2).Value = "Hospital San Jorge"
$ws.Cells.Item(2931, 3).Value = 27
$ws.Cells.Item(2931, 4).Value = 7
$ws.Cells.Item(2931, 5).Value = "Huesca"
$ws.Cells.Item(2931, 6).Value = "Huesca"
$ws.Cells.Item(2931, 7).Value = 22125
$ws.Cells.Item(2931, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2932, 1).Value = 44080
$ws.Cells.Item(2932, 2).Value = "Hospital Sagrado Corazón"
$ws.Cells.Item(2932, 3).Value = 0
$ws.Cells.Item(2932, 4).Value = 0
$ws.Cells.Item(2932, 5).Value = "Huesca"
$ws.Cells.Item(2932, 6).Value = "Huesca"
$ws.Cells.Item(2932, 7).Value = 22125
$ws.Cells.Item(2932, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2933, 1).Value = 44080
$ws.Cells.Item(2933, 2).Value = "Hospital Ernest Lluch"
$ws.Cells.Item(2933, 3).Value = 4
$ws.Cells.Item(2933, 4).Value = 0
$ws.Cells.Item(2933, 5).Value = "Calatayud"
$ws.Cells.Item(2933, 6).Value = "Zaragoza"
$ws.Cells.Item(2933, 7).Value = 50297
$ws.Cells.Item(2933, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2934, 1).Value = 44080
$ws.Cells.Item(2934, 2).Value = "Hospital San José"
$ws.Cells.Item(2934, 3).Value = 1
$ws.Cells.Item(2934, 4).Value = 0
$ws.Cells.Item(2934, 5).Value = "Teruel"
$ws.Cells.Item(2934, 6).Value = "Teruel"
$ws.Cells.Item(2934, 7).Value = 44216
$ws.Cells.Item(2934, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2935, 1).Value = 44080
$ws.Cells.Item(2935, 2).Value = "Hospital Ejea – Cinco Villas"
$ws.Cells.Item(2935, 3).Value = 1
$ws.Cells.Item(2935, 4).Value = 0
$ws.Cells.Item(2935, 5).Value = "Ejea de los Caballeros"
$ws.Cells.Item(2935, 6).Value = "Zaragoza"
$ws.Cells.Item(2935, 7).Value = 50297
$ws.Cells.Item(2935, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2936, 1).Value = 44080
$ws.Cells.Item(2936, 2).Value = "MAZ"
$ws.Cells.Item(2936, 3).Value = 2
$ws.Cells.Item(2936, 4).Value = 0
$ws.Cells.Item(2936, 5).Value = "Zaragoza"
$ws.Cells.Item(2936, 6).Value = "Zaragoza"
$ws.Cells.Item(2936, 7).Value = 50297
$ws.Cells.Item(2936, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2937, 1).Value = 44080
$ws.Cells.Item(2937, 2).Value = "Hospital Viamed Montecanal"
$ws.Cells.Item(2937, 3).Value = 0
$ws.Cells.Item(2937, 4).Value = 0
$ws.Cells.Item(2937, 5).Value = "Zaragoza"
$ws.Cells.Item(2937, 6).Value = "Zaragoza"
$ws.Cells.Item(2937, 7).Value = 50297
$ws.Cells.Item(2937, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2938, 1).Value = 44080
$ws.Cells.Item(2938, 2).Value = "Clínica Montpellier"
$ws.Cells.Item(2938, 3).Value = 4
$ws.Cells.Item(2938, 4).Value = 0
$ws.Cells.Item(2938, 5).Value = "Zaragoza"
$ws.Cells.Item(2938, 6).Value = "Zaragoza"
$ws.Cells.Item(2938, 7).Value = 50297
$ws.Cells.Item(2938, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2939, 1).Value = 44080
$ws.Cells.Item(2939, 2).Value = "Hospital Quirón"
$ws.Cells.Item(2939, 3).Value = 6
$ws.Cells.Item(2939, 4).Value = 2
$ws.Cells.Item(2939, 5).Value = "Zaragoza"
$ws.Cells.Item(2939, 6).Value = "Zaragoza"
$ws.Cells.Item(2939, 7).Value = 50297
$ws.Cells.Item(2939, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2940, 1).Value = 44080
$ws.Cells.Item(2940, 2).Value = "Hospital San Juan de Dios de Zaragoza"
$ws.Cells.Item(2940, 3).Value = 35
$ws.Cells.Item(2940, 4).Value = 0
$ws.Cells.Item(2940, 5).Value = "Zaragoza"
$ws.Cells.Item(2940, 6).Value = "Zaragoza"
$ws.Cells.Item(2940, 7).Value = 50297
$ws.Cells.Item(2940, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2941, 1).Value = 44080
$ws.Cells.Item(2941, 2).Value = "Clínica Viamed Santiago"
$ws.Cells.Item(2941, 3).Value = 2
$ws.Cells.Item(2941, 4).Value = 0
$ws.Cells.Item(2941, 5).Value = "Huesca"
$ws.Cells.Item(2941, 6).Value = "Huesca"
$ws.Cells.Item(2941, 7).Value = 22125
$ws.Cells.Item(2941, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2942, 1).Value = 44080
$ws.Cells.Item(2942, 2).Value = "Clínica El Pilar"
$ws.Cells.Item(2942, 3).Value = 1
$ws.Cells.Item(2942, 4).Value = 0
$ws.Cells.Item(2942, 5).Value = "Zaragoza"
$ws.Cells.Item(2942, 6).Value = "Zaragoza"
$ws.Cells.Item(2942, 7).Value = 50297
$ws.Cells.Item(2942, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2943, 1).Value = 44080
$ws.Cells.Item(2943, 2).Value = "C.S. Fraga - Bajo Cinca"
$ws.Cells.Item(2943, 3).Value = 8
$ws.Cells.Item(2943, 4).Value = 0
$ws.Cells.Item(2943, 5).Value = "Fraga"
$ws.Cells.Item(2943, 6).Value = "Huesca"
$ws.Cells.Item(2943, 7).Value = 22125
$ws.Cells.Item(2943, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2944, 1).Value = 44080
$ws.Cells.Item(2944, 2).Value = "Hospital de Jaca"
$ws.Cells.Item(2944, 3).Value = 0
$ws.Cells.Item(2944, 4).Value = 0
$ws.Cells.Item(2944, 5).Value = "Jaca"
$ws.Cells.Item(2944, 6).Value = "Huesca"
$ws.Cells.Item(2944, 7).Value = 22125
$ws.Cells.Item(2944, 8).Value = "Fuente Servcio Aragonés de Salud - 20200905 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2945, 1).Value = 44081
$ws.Cells.Item(2945, 2).Value = "Hospital Universitario Miguel Servet"
$ws.Cells.Item(2945, 3).Value = 109
$ws.Cells.Item(2945, 4).Value = 23
$ws.Cells.Item(2945, 5).Value = "Zaragoza"
$ws.Cells.Item(2945, 6).Value = "Zaragoza"
$ws.Cells.Item(2945, 7).Value = 50297
$ws.Cells.Item(2945, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2946, 1).Value = 44081
$ws.Cells.Item(2946, 2).Value = "Hospital Clínico Universitario"
$ws.Cells.Item(2946, 3).Value = 158
$ws.Cells.Item(2946, 4).Value = 23
$ws.Cells.Item(2946, 5).Value = "Zaragoza"
$ws.Cells.Item(2946, 6).Value = "Zaragoza"
$ws.Cells.Item(2946, 7).Value = 50297
$ws.Cells.Item(2946, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2947, 1).Value = 44081
$ws.Cells.Item(2947, 2).Value = "Hospital Royo Villanova"
$ws.Cells.Item(2947, 3).Value = 46
$ws.Cells.Item(2947, 4).Value = 4
$ws.Cells.Item(2947, 5).Value = "Zaragoza"
$ws.Cells.Item(2947, 6).Value = "Zaragoza"
$ws.Cells.Item(2947, 7).Value = 50297
$ws.Cells.Item(2947, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2948, 1).Value = 44081
$ws.Cells.Item(2948, 2).Value = "Hospital Nuestra Señora de Gracia"
$ws.Cells.Item(2948, 3).Value = 12
$ws.Cells.Item(2948, 4).Value = 1
$ws.Cells.Item(2948, 5).Value = "Zaragoza"
$ws.Cells.Item(2948, 6).Value = "Zaragoza"
$ws.Cells.Item(2948, 7).Value = 50297
$ws.Cells.Item(2948, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2949, 1).Value = 44081
$ws.Cells.Item(2949, 2).Value = "Hospital General de la Defensa"
$ws.Cells.Item(2949, 3).Value = 16
$ws.Cells.Item(2949, 4).Value = 2
$ws.Cells.Item(2949, 5).Value = "Zaragoza"
$ws.Cells.Item(2949, 6).Value = "Zaragoza"
$ws.Cells.Item(2949, 7).Value = 50297
$ws.Cells.Item(2949, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2950, 1).Value = 44081
$ws.Cells.Item(2950, 2).Value = "Hospital Obispo Polanco"
$ws.Cells.Item(2950, 3).Value = 15
$ws.Cells.Item(2950, 4).Value = 1
$ws.Cells.Item(2950, 5).Value = "Teruel"
$ws.Cells.Item(2950, 6).Value = "Teruel"
$ws.Cells.Item(2950, 7).Value = 44216
$ws.Cells.Item(2950, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2951, 1).Value = 44081
$ws.Cells.Item(2951, 2).Value = "Hospital de Alcañiz"
$ws.Cells.Item(2951, 3).Value = 23
$ws.Cells.Item(2951, 4).Value = 0
$ws.Cells.Item(2951, 5).Value = "Alcañiz"
$ws.Cells.Item(2951, 6).Value = "Teruel"
$ws.Cells.Item(2951, 7).Value = 44216
$ws.Cells.Item(2951, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2952, 1).Value = 44081
$ws.Cells.Item(2952, 2).Value = "Hospital de Barbastro"
$ws.Cells.Item(2952, 3).Value = 13
$ws.Cells.Item(2952, 4).Value = 3
$ws.Cells.Item(2952, 5).Value = "Barbastro"
$ws.Cells.Item(2952, 6).Value = "Huesca"
$ws.Cells.Item(2952, 7).Value = 22125
$ws.Cells.Item(2952, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2953, 1).Value = 44081
$ws.Cells.Item(2953, 2).Value = "Hospital San Jorge"
$ws.Cells.Item(2953, 3).Value = 28
$ws.Cells.Item(2953, 4).Value = 7
$ws.Cells.Item(2953, 5).Value = "Huesca"
$ws.Cells.Item(2953, 6).Value = "Huesca"
$ws.Cells.Item(2953, 7).Value = 22125
$ws.Cells.Item(2953, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2954, 1).Value = 44081
$ws.Cells.Item(2954, 2).Value = "Hospital Sagrado Corazón"
$ws.Cells.Item(2954, 3).Value = 0
$ws.Cells.Item(2954, 4).Value = 0
$ws.Cells.Item(2954, 5).Value = "Huesca"
$ws.Cells.Item(2954, 6).Value = "Huesca"
$ws.Cells.Item(2954, 7).Value = 22125
$ws.Cells.Item(2954, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2955, 1).Value = 44081
$ws.Cells.Item(2955, 2).Value = "Hospital Ernest Lluch"
$ws.Cells.Item(2955, 3).Value = 4
$ws.Cells.Item(2955, 4).Value = 0
$ws.Cells.Item(2955, 5).Value = "Calatayud"
$ws.Cells.Item(2955, 6).Value = "Zaragoza"
$ws.Cells.Item(2955, 7).Value = 50297
$ws.Cells.Item(2955, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2956, 1).Value = 44081
$ws.Cells.Item(2956, 2).Value = "Hospital San José"
$ws.Cells.Item(2956, 3).Value = 1
$ws.Cells.Item(2956, 4).Value = 0
$ws.Cells.Item(2956, 5).Value = "Teruel"
$ws.Cells.Item(2956, 6).Value = "Teruel"
$ws.Cells.Item(2956, 7).Value = 44216
$ws.Cells.Item(2956, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2957, 1).Value = 44081
$ws.Cells.Item(2957, 2).Value = "Hospital Ejea – Cinco Villas"
$ws.Cells.Item(2957, 3).Value = 1
$ws.Cells.Item(2957, 4).Value = 0
$ws.Cells.Item(2957, 5).Value = "Ejea de los Caballeros"
$ws.Cells.Item(2957, 6).Value = "Zaragoza"
$ws.Cells.Item(2957, 7).Value = 50297
$ws.Cells.Item(2957, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2958, 1).Value = 44081
$ws.Cells.Item(2958, 2).Value = "MAZ"
$ws.Cells.Item(2958, 3).Value = 2
$ws.Cells.Item(2958, 4).Value = 0
$ws.Cells.Item(2958, 5).Value = "Zaragoza"
$ws.Cells.Item(2958, 6).Value = "Zaragoza"
$ws.Cells.Item(2958, 7).Value = 50297
$ws.Cells.Item(2958, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2959, 1).Value = 44081
$ws.Cells.Item(2959, 2).Value = "Hospital Viamed Montecanal"
$ws.Cells.Item(2959, 3).Value = 0
$ws.Cells.Item(2959, 4).Value = 0
$ws.Cells.Item(2959, 5).Value = "Zaragoza"
$ws.Cells.Item(2959, 6).Value = "Zaragoza"
$ws.Cells.Item(2959, 7).Value = 50297
$ws.Cells.Item(2959, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2960, 1).Value = 44081
$ws.Cells.Item(2960, 2).Value = "Clínica Montpellier"
$ws.Cells.Item(2960, 3).Value = 4
$ws.Cells.Item(2960, 4).Value = 0
$ws.Cells.Item(2960, 5).Value = "Zaragoza"
$ws.Cells.Item(2960, 6).Value = "Zaragoza"
$ws.Cells.Item(2960, 7).Value = 50297
$ws.Cells.Item(2960, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2961, 1).Value = 44081
$ws.Cells.Item(2961, 2).Value = "Hospital Quirón"
$ws.Cells.Item(2961, 3).Value = 7
$ws.Cells.Item(2961, 4).Value = 2
$ws.Cells.Item(2961, 5).Value = "Zaragoza"
$ws.Cells.Item(2961, 6).Value = "Zaragoza"
$ws.Cells.Item(2961, 7).Value = 50297
$ws.Cells.Item(2961, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2962, 1).Value = 44081
$ws.Cells.Item(2962, 2).Value = "Hospital San Juan de Dios de Zaragoza"
$ws.Cells.Item(2962, 3).Value = 35
$ws.Cells.Item(2962, 4).Value = 0
$ws.Cells.Item(2962, 5).Value = "Zaragoza"
$ws.Cells.Item(2962, 6).Value = "Zaragoza"
$ws.Cells.Item(2962, 7).Value = 50297
$ws.Cells.Item(2962, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2963, 1).Value = 44081
$ws.Cells.Item(2963, 2).Value = "Clínica Viamed Santiago"
$ws.Cells.Item(2963, 3).Value = 2
$ws.Cells.Item(2963, 4).Value = 0
$ws.Cells.Item(2963, 5).Value = "Huesca"
$ws.Cells.Item(2963, 6).Value = "Huesca"
$ws.Cells.Item(2963, 7).Value = 22125
$ws.Cells.Item(2963, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2964, 1).Value = 44081
$ws.Cells.Item(2964, 2).Value = "Clínica El Pilar"
$ws.Cells.Item(2964, 3).Value = 1
$ws.Cells.Item(2964, 4).Value = 0
$ws.Cells.Item(2964, 5).Value = "Zaragoza"
$ws.Cells.Item(2964, 6).Value = "Zaragoza"
$ws.Cells.Item(2964, 7).Value = 50297
$ws.Cells.Item(2964, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2965, 1).Value = 44081
$ws.Cells.Item(2965, 2).Value = "C.S. Fraga - Bajo Cinca"
$ws.Cells.Item(2965, 3).Value = 8
$ws.Cells.Item(2965, 4).Value = 0
$ws.Cells.Item(2965, 5).Value = "Fraga"
$ws.Cells.Item(2965, 6).Value = "Huesca"
$ws.Cells.Item(2965, 7).Value = 22125
$ws.Cells.Item(2965, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

$ws.Cells.Item(2966, 1).Value = 44081
$ws.Cells.Item(2966, 2).Value = "Hospital de Jaca"
$ws.Cells.Item(2966, 3).Value = 0
$ws.Cells.Item(2966, 4).Value = 0
$ws.Cells.Item(2966, 5).Value = "Jaca"
$ws.Cells.Item(2966, 6).Value = "Huesca"
$ws.Cells.Item(2966, 7).Value = 22125
$ws.Cells.Item(2966, 8).Value = "Fuente Servcio Aragonés de Salud - 20200906 COVID-19 SITUACION HOSPITALIZACION.xlsx"

# Update the selection / scroll position to mirror the author's last view
$null = $ws.Range("A2946:A2966").Select()
